# RPAR_holdings.xlsx update:
#  - Refresh the "as of" date in the confidential disclaimer text (A18)
#  - Refresh Weight (col D) and Percent Change (col E) figures for rows 2-15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so the cells can be written, then
# restore protection afterwards.
$ws.Unprotect()

# --- Update the disclaimer "as of" date -------------------------------------
$disclaimer = $ws.Range("A18").Value()
$disclaimer = $disclaimer -replace "2021-06-14", "2021-07-07"
$ws.Range("A18").Value = $disclaimer

# --- Update Weight (D) / Percent Change (E) values --------------------------
$ws.Range("D2").Value = 0.05872398512175927
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 0.02020122853252084
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.02768737406888534
$ws.Range("E4").Value = 0

$ws.Range("D5").Value = 0.02914580741234215
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.02908376031532988
$ws.Range("E6").Value = 0

$ws.Range("D7").Value = 0.01800868006230835
$ws.Range("E7").Value = 0

$ws.Range("D8").Value = 0.01081252306355909
$ws.Range("E8").Value = 0

$ws.Range("D9").Value = 0.01059437853301069
$ws.Range("E9").Value = 0

$ws.Range("D10").Value = 0.0664524409001401
$ws.Range("E10").Value = 0

$ws.Range("D11").Value = 0.06656412567476218
$ws.Range("E11").Value = 0

$ws.Range("D12").Value = 0.1533438486834586
$ws.Range("E12").Value = 0

$ws.Range("D13").Value = 0.3951890640358698
$ws.Range("E13").Value = 0

$ws.Range("D14").Value = 0.1141927835960538
$ws.Range("E14").Value = 0

$ws.Range("E15").Value = 0

# Restore sheet protection to match the original workbook state.
$ws.Protect()
